$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff swaps the species-related data (columns A, E, F, G, H, Q, R) between
# row 2 and row 3, while column B gets new distinct values for each row.

# --- Row 2 new values ---
$ws.Range("A2").Value = 111950173
$ws.Range("B2").Value = 90806
$ws.Range("E2").Value = 4361
$ws.Range("F2").Value = "Orange taggsvamp"
$ws.Range("G2").Value = "Hydnellum aurantiacum"
$ws.Range("H2").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q2").Value = 465440
$ws.Range("R2").Value = 6875680

# --- Row 3 new values ---
$ws.Range("A3").Value = 111950243
$ws.Range("B3").Value = 90837
$ws.Range("E3").Value = 5966
$ws.Range("F3").Value = "Motaggsvamp"
$ws.Range("G3").Value = "Sarcodon squamosus"
$ws.Range("H3").Value = "(Schaeff.) Quél."
$ws.Range("Q3").Value = 465473
$ws.Range("R3").Value = 6875785
